$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the existing last column (admin_comment),
# pushing it from O to Q, and creating new O/P columns for the
# business manager fields.
$ws.Columns("O:P").Insert()

# New header cells
$ws.Range("O1").Value = "business_manager_name"
$ws.Range("P1").Value = "business_manager_mobile_no"

# New data cells for row 2
$ws.Range("O2").Value = "MMM"
$ws.Range("P2").Value = 9835677898

# Update selection to reflect the new active cell
$ws.Range("O2").Select()
